$d = $word.ActiveDocument

# The bibliography entry for Bertero is followed by three paragraphs that
# must be removed in their entirety:
#   1) an empty "Normal" paragraph
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#      pages. Original theme under Creative Commons Attribution"
# Everything before and after this block (including the empty paragraph
# and page-break paragraph that follow it) must stay untouched.

$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text

    if ($t -match "Ver no Jupiter Salvar em pdf Salvar em docx") {
        $startPara = $d.Paragraphs.Item($i - 1)
    }
    if ($t -match "Contact: luizeleno@usp\.br") {
        $endPara = $p
    }
}

$rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
$rng.Delete()
